$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 14 (existing rows 14-41 shift down to 16-43).
$ws.Rows.Item(14).Resize(2).Insert()

# New row 14 data
$ws.Cells.Item(14, 1).Value = 1
$ws.Cells.Item(14, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(14, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(14, 4).Value = 44868
$ws.Cells.Item(14, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(14, 5).Value = 15
$ws.Cells.Item(14, 6).Value = 100112028
$ws.Cells.Item(14, 7).Value = "Sandia"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 600
$ws.Cells.Item(14, 12).Value = 630
$ws.Cells.Item(14, 13).Value = 615
$ws.Cells.Item(14, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(14, 15).Value = "Perú"
$ws.Cells.Item(14, 16).Value = 615
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"

# New row 15 data
$ws.Cells.Item(15, 1).Value = 1
$ws.Cells.Item(15, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(15, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(15, 4).Value = 44868
$ws.Cells.Item(15, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(15, 5).Value = 15
$ws.Cells.Item(15, 6).Value = 100112028
$ws.Cells.Item(15, 7).Value = "Sandia"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Segunda"
$ws.Cells.Item(15, 10).Value = 300
$ws.Cells.Item(15, 11).Value = 580
$ws.Cells.Item(15, 12).Value = 600
$ws.Cells.Item(15, 13).Value = 590
$ws.Cells.Item(15, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 590
$ws.Cells.Item(15, 17).Value = 1
$ws.Cells.Item(15, 18).Value = "Hortaliza"

# Update the sheet dimension to reflect the new used range.
$ws.Range("A1:R43").Select()
